# Apply "wordpress version - capitulo 8 esp" update:
# Set status (column C) to "wordpress esp" for row 28 (Chapter 7 - Potencia)
# and row 29 (Chapter 8 - Costes), and move the active selection to C29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C28").Value = "wordpress esp"
$ws.Range("C29").Value = "wordpress esp"

$ws.Range("C29").Select()
